$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.127.28'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '1.850.79'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4649'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2804'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06404'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '96.41'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +13.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07539'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = '1.841.70'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.969'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6335'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '293.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +21.22%  '
$ws.Range('D17').Value = '30.101.95'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.003'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007359'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.005'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.076.67'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.994'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.037'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.073'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.922'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1080'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.335'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.003'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.807'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04909'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7239'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.108'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.734'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01917'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.654'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8634'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.959'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.003'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.607'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4034'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.026'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.960'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1184'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05551'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '840.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +18.27%  '
